$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.2
$ws.Range("H2").Value = 3.4
$ws.Range("J2").Value = 2.88
$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 2.02
$ws.Range("U2").Value = 1.36
$ws.Range("V2").Value = 3
$ws.Range("W2").Value = 1.67
$ws.Range("X2").Value = 2.1
$ws.Range("AD2").Value = 23
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 6.5
$ws.Range("AI2").Value = 151
$ws.Range("AO2").Value = 29
# Row 3
$ws.Range("G3").Value = 2.3
$ws.Range("I3").Value = 3.5
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 4
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.75
$ws.Range("W3").Value = 2
$ws.Range("X3").Value = 1.75
$ws.Range("Y3").Value = 6.5
$ws.Range("AI3").Value = 451
$ws.Range("AP3").Value = 1.87
$ws.Range("AQ3").Value = 2.03
# Row 4
$ws.Range("G4").Value = 2.38
$ws.Range("I4").Value = 3.25
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 3.75
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.62
$ws.Range("S4").Value = 4.33
$ws.Range("T4").Value = 1.22
$ws.Range("W4").Value = 1.95
$ws.Range("X4").Value = 1.8
$ws.Range("Y4").Value = 7
$ws.Range("AE4").Value = 7.5
$ws.Range("AI4").Value = 351
$ws.Range("AN4").Value = 29
# Row 5
$ws.Range("G5").Value = 2.5
$ws.Range("K5").Value = 1.83
$ws.Range("L5").Value = 4.33
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("Y5").Value = 5.5
$ws.Range("AB5").Value = 23
$ws.Range("AK5").Value = 15
$ws.Range("AM5").Value = 41
$ws.Range("AP5").Value = 2.14
$ws.Range("AR5").Value = 4.93
# Row 6
$ws.Range("G6").Value = 2.05
$ws.Range("J6").Value = 2.88
$ws.Range("K6").Value = 1.91
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("O6").Value = 1.5
$ws.Range("P6").Value = 2.5
$ws.Range("Q6").Value = 2.6
$ws.Range("R6").Value = 1.48
$ws.Range("S6").Value = 5.5
$ws.Range("T6").Value = 1.14
$ws.Range("W6").Value = 2.2
$ws.Range("X6").Value = 1.62
$ws.Range("AA6").Value = 10
$ws.Range("AC6").Value = 21
$ws.Range("AE6").Value = 6
$ws.Range("AH6").Value = 81
$ws.Range("AJ6").Value = 9
$ws.Range("AP6").Value = 2
$ws.Range("AQ6").Value = 1.85
$ws.Range("AR6").Value = 4.1
$ws.Range("AS6").Value = 1.23
# Row 7
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
# Row 11
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 7
$ws.Range("AP11").Value = 1.8
$ws.Range("AQ11").Value = 2.05
# Row 12
$ws.Range("G12").Value = 2.88
$ws.Range("H12").Value = 2.75
$ws.Range("I12").Value = 2.45
$ws.Range("K12").Value = 1.95
$ws.Range("L12").Value = 3.25
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.63
$ws.Range("Q12").Value = 2.4
$ws.Range("R12").Value = 1.53
$ws.Range("S12").Value = 4.5
$ws.Range("T12").Value = 1.18
$ws.Range("U12").Value = 1.53
$ws.Range("V12").Value = 2.38
$ws.Range("W12").Value = 2
$ws.Range("X12").Value = 1.73
$ws.Range("Y12").Value = 8
$ws.Range("Z12").Value = 13
$ws.Range("AB12").Value = 29
$ws.Range("AE12").Value = 7
$ws.Range("AL12").Value = 11
$ws.Range("AM12").Value = 26
$ws.Range("AN12").Value = 23
$ws.Range("AO12").Value = 41
$ws.Range("AP12").Value = 1.8
$ws.Range("AQ12").Value = 2.05
# Row 15
$ws.Range("K15").Value = 2.4
$ws.Range("Q15").Value = 1.53
$ws.Range("R15").Value = 2.4
$ws.Range("S15").Value = 2.25
$ws.Range("T15").Value = 1.57
$ws.Range("W15").Value = 1.44
$ws.Range("X15").Value = 2.63
$ws.Range("Z15").Value = 15
$ws.Range("AE15").Value = 19
$ws.Range("AG15").Value = 11
$ws.Range("AO15").Value = 21
$ws.Range("AR15").Value = 1.88
$ws.Range("AS15").Value = 1.98
# Row 16
$ws.Range("G16").Value = 1.55
$ws.Range("H16").Value = 4.33
$ws.Range("J16").Value = 2.05
$ws.Range("K16").Value = 2.5
$ws.Range("L16").Value = 5.5
$ws.Range("N16").Value = 17
$ws.Range("Q16").Value = 1.57
$ws.Range("R16").Value = 2.35
$ws.Range("AF16").Value = 8.5
$ws.Range("AR16").Value = 2
$ws.Range("AS16").Value = 1.85
# Row 17
$ws.Range("G17").Value = 2.3
$ws.Range("H17").Value = 3.9
$ws.Range("I17").Value = 2.7
$ws.Range("L17").Value = 3.1
$ws.Range("Y17").Value = 15
$ws.Range("AA17").Value = 10
$ws.Range("AJ17").Value = 15
# Row 18
$ws.Range("G18").Value = 3
$ws.Range("I18").Value = 2.5
$ws.Range("L18").Value = 3.25
$ws.Range("O18").Value = 1.4
$ws.Range("P18").Value = 2.75
$ws.Range("W18").Value = 1.95
$ws.Range("X18").Value = 1.8
$ws.Range("Z18").Value = 13
$ws.Range("AC18").Value = 26
$ws.Range("AN18").Value = 23
# Row 19
$ws.Range("G19").Value = 2.35
$ws.Range("I19").Value = 2.8
$ws.Range("R19").Value = 1.75
$ws.Range("S19").Value = 3.75
$ws.Range("T19").Value = 1.25
$ws.Range("W19").Value = 1.83
$ws.Range("X19").Value = 1.83
$ws.Range("AB19").Value = 23
$ws.Range("AC19").Value = 21
$ws.Range("AK19").Value = 13
# Row 23
$ws.Range("G23").Value = 2.05
$ws.Range("I23").Value = 4
$ws.Range("L23").Value = 4.75
$ws.Range("O23").Value = 1.53
$ws.Range("P23").Value = 2.38
$ws.Range("Q23").Value = 2.7
$ws.Range("R23").Value = 1.44
$ws.Range("S23").Value = 5.5
$ws.Range("T23").Value = 1.14
$ws.Range("W23").Value = 2.25
$ws.Range("X23").Value = 1.57
$ws.Range("Y23").Value = 5.5
$ws.Range("Z23").Value = 8.5
$ws.Range("AB23").Value = 19
$ws.Range("AL23").Value = 15
$ws.Range("AM23").Value = 41
$ws.Range("AP23").Value = 2.05
$ws.Range("AQ23").Value = 1.8
